$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the last row (row 53), shifting dimension to A1:E52
$ws.Rows.Item(53).Delete()

# Update data rows 2-52 with corrected values
$ws.Range("A2").Value = 39583
$ws.Range("B2").Value = 2008
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 2009
$ws.Range("E2").Value = 4.036825632626817

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = 3.407109591918855

$ws.Range("A4").Value = 39948
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -3.591129714716879
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = -1.300150869559236

$ws.Range("A5").Value = 40130
$ws.Range("B5").Value = 2009
$ws.Range("C5").ClearContents()
$ws.Range("D5").Value = 2010
$ws.Range("E5").Value = -2.725947775269033

$ws.Range("A6").Value = 40310
$ws.Range("B6").Value = 2010
$ws.Range("C6").Value = 3.52827217675542
$ws.Range("D6").Value = 2011
$ws.Range("E6").Value = 1.194245528719495

$ws.Range("A7").Value = 40494
$ws.Range("B7").Value = 2010
$ws.Range("C7").ClearContents()
$ws.Range("D7").Value = 2011
$ws.Range("E7").Value = 5.993806847197725

$ws.Range("A8").Value = 40676
$ws.Range("B8").Value = 2011
$ws.Range("C8").Value = 7.607887362976751
$ws.Range("D8").Value = 2012
$ws.Range("E8").Value = 3.063639588842682

$ws.Range("A9").Value = 40862
$ws.Range("B9").Value = 2011
$ws.Range("C9").ClearContents()
$ws.Range("D9").Value = 2012
$ws.Range("E9").Value = 6.924353497010971

$ws.Range("A10").Value = 41044
$ws.Range("B10").Value = 2012
$ws.Range("C10").Value = 3.881953143326466
$ws.Range("D10").Value = 2013
$ws.Range("E10").Value = 7.155859501153827

$ws.Range("A11").Value = 41228
$ws.Range("B11").Value = 2012
$ws.Range("C11").ClearContents()
$ws.Range("D11").Value = 2013
$ws.Range("E11").Value = 3.737237443362851

$ws.Range("A12").Value = 41409
$ws.Range("B12").Value = 2013
$ws.Range("C12").Value = 0.354775830825127
$ws.Range("D12").Value = 2014
$ws.Range("E12").Value = 2.700497159199755

$ws.Range("A13").Value = 41592
$ws.Range("B13").Value = 2013
$ws.Range("C13").ClearContents()
$ws.Range("D13").Value = 2014
$ws.Range("E13").Value = 3.504647804006344

$ws.Range("A14").Value = 41774
$ws.Range("B14").Value = 2014
$ws.Range("C14").Value = 5.106323395421475
$ws.Range("D14").Value = 2015
$ws.Range("E14").Value = 2.62739064366051

$ws.Range("A15").Value = 41957
$ws.Range("B15").Value = 2014
$ws.Range("C15").ClearContents()
$ws.Range("D15").Value = 2015
$ws.Range("E15").Value = 4.130094879572455

$ws.Range("A16").Value = 42137
$ws.Range("B16").Value = 2015
$ws.Range("C16").Value = 3.830515520137801
$ws.Range("D16").Value = 2016
$ws.Range("E16").Value = 3.858663141671226

$ws.Range("A17").Value = 42321
$ws.Range("B17").Value = 2015
$ws.Range("C17").ClearContents()
$ws.Range("D17").Value = 2016
$ws.Range("E17").Value = 5.194458387461709

$ws.Range("A18").Value = 42503
$ws.Range("B18").Value = 2016
$ws.Range("C18").Value = 4.328608026086478
$ws.Range("D18").Value = 2017
$ws.Range("E18").Value = 4.721407739775696

$ws.Range("A19").Value = 42689
$ws.Range("B19").Value = 2016
$ws.Range("C19").ClearContents()
$ws.Range("D19").Value = 2017
$ws.Range("E19").Value = 4.372655645302403

$ws.Range("A20").Value = 42867
$ws.Range("B20").Value = 2017
$ws.Range("C20").Value = 4.429102498614346
$ws.Range("D20").Value = 2018
$ws.Range("E20").Value = 4.21218881008929

$ws.Range("A21").Value = 43053
$ws.Range("B21").Value = 2017
$ws.Range("C21").ClearContents()
$ws.Range("D21").Value = 2018
$ws.Range("E21").Value = 5.324897060120137

$ws.Range("A22").Value = 43145
$ws.Range("B22").Value = 2018
$ws.Range("C22").Value = 5.783465271898192
$ws.Range("D22").Value = 2019
$ws.Range("E22").Value = 4.964333730716497

$ws.Range("A23").Value = 43235
$ws.Range("B23").Value = 2018
$ws.Range("C23").Value = 5.850954342715009
$ws.Range("D23").Value = 2019
$ws.Range("E23").Value = 4.990046926794744

$ws.Range("A24").Value = 43326
$ws.Range("B24").Value = 2018
$ws.Range("C24").Value = 5.410900500218596
$ws.Range("D24").Value = 2019
$ws.Range("E24").Value = 4.531726208768672

$ws.Range("A25").Value = 43418
$ws.Range("B25").Value = 2018
$ws.Range("C25").ClearContents()
$ws.Range("D25").Value = 2019
$ws.Range("E25").Value = 4.693063499664252

$ws.Range("A26").Value = 43510
$ws.Range("B26").Value = 2019
$ws.Range("C26").Value = 4.115125864415514
$ws.Range("D26").Value = 2020
$ws.Range("E26").Value = 4.871235017471043

$ws.Range("A27").Value = 43600
$ws.Range("B27").Value = 2019
$ws.Range("C27").Value = 4.557673974453769
$ws.Range("D27").Value = 2020
$ws.Range("E27").Value = 5.208344373007368

$ws.Range("A28").Value = 43691
$ws.Range("B28").Value = 2019
$ws.Range("C28").Value = 3.383644094252025
$ws.Range("D28").Value = 2020
$ws.Range("E28").Value = 3.4540322565767

$ws.Range("A29").Value = 43783
$ws.Range("B29").Value = 2019
$ws.Range("C29").ClearContents()
$ws.Range("D29").Value = 2020
$ws.Range("E29").Value = 3.273620772016161

$ws.Range("A30").Value = 43875
$ws.Range("B30").Value = 2020
$ws.Range("C30").Value = 3.2928463695165
$ws.Range("D30").Value = 2021
$ws.Range("E30").Value = 3.901680317162204

$ws.Range("A31").Value = 43966
$ws.Range("B31").Value = 2020
$ws.Range("C31").Value = 1.13158575217045
$ws.Range("D31").Value = 2021
$ws.Range("E31").Value = 2.143123507515932

$ws.Range("A32").Value = 44068
$ws.Range("B32").Value = 2020
$ws.Range("C32").Value = -9.2489161297999
$ws.Range("D32").Value = 2021
$ws.Range("E32").Value = -8.033751002286738

$ws.Range("A33").Value = 44159
$ws.Range("B33").Value = 2020
$ws.Range("C33").ClearContents()
$ws.Range("D33").Value = 2021
$ws.Range("E33").Value = -4.319815935184923

$ws.Range("A34").Value = 44251
$ws.Range("B34").Value = 2021
$ws.Range("C34").Value = -2.573749071621145
$ws.Range("D34").Value = 2022
$ws.Range("E34").Value = -1.545753571259545

$ws.Range("A35").Value = 44341
$ws.Range("B35").Value = 2021
$ws.Range("C35").Value = -1.77012120409461
$ws.Range("D35").Value = 2022
$ws.Range("E35").Value = -0.6252235182164778

$ws.Range("A36").Value = 44432
$ws.Range("B36").Value = 2021
$ws.Range("C36").Value = -1.287084480507283
$ws.Range("D36").Value = 2022
$ws.Range("E36").Value = 0.7555897036957804

$ws.Range("A37").Value = 44525
$ws.Range("B37").Value = 2021
$ws.Range("C37").ClearContents()
$ws.Range("D37").Value = 2022
$ws.Range("E37").Value = 1.850145327219943

$ws.Range("A38").Value = 44617
$ws.Range("B38").Value = 2022
$ws.Range("C38").Value = 0.3839706909697815
$ws.Range("D38").Value = 2023
$ws.Range("E38").Value = -3.055292064789206

$ws.Range("A39").Value = 44706
$ws.Range("B39").Value = 2022
$ws.Range("C39").Value = 1.286283684448075
$ws.Range("D39").Value = 2023
$ws.Range("E39").Value = -2.183761975384579

$ws.Range("A40").Value = 44798
$ws.Range("B40").Value = 2022
$ws.Range("C40").Value = 1.494343500592232
$ws.Range("D40").Value = 2023
$ws.Range("E40").Value = -1.161086276452006

$ws.Range("A41").Value = 44890
$ws.Range("B41").Value = 2022
$ws.Range("C41").ClearContents()
$ws.Range("D41").Value = 2023
$ws.Range("E41").Value = 2.552834403233084

$ws.Range("A42").Value = 44981
$ws.Range("B42").Value = 2023
$ws.Range("C42").Value = -0.3407161957438287
$ws.Range("D42").Value = 2024
$ws.Range("E42").Value = 1.987749514178372

$ws.Range("A43").Value = 45071
$ws.Range("B43").Value = 2023
$ws.Range("C43").Value = -0.7189954590872905
$ws.Range("D43").Value = 2024
$ws.Range("E43").Value = 0.7347074005453758

$ws.Range("A44").Value = 45163
$ws.Range("B44").Value = 2023
$ws.Range("C44").Value = -0.6982718287330991
$ws.Range("D44").Value = 2024
$ws.Range("E44").Value = 0.4357420177721227

$ws.Range("A45").Value = 45254
$ws.Range("B45").Value = 2023
$ws.Range("C45").ClearContents()
$ws.Range("D45").Value = 2024
$ws.Range("E45").Value = -1.04392885455985

$ws.Range("A46").Value = 45345
$ws.Range("B46").Value = 2024
$ws.Range("C46").Value = -1.305399582732825
$ws.Range("D46").Value = 2025
$ws.Range("E46").Value = -0.9722577196979332

$ws.Range("A47").Value = 45436
$ws.Range("B47").Value = 2024
$ws.Range("C47").Value = -0.382605475081077
$ws.Range("D47").Value = 2025
$ws.Range("E47").Value = 0.3097078768351302

$ws.Range("A48").Value = 45534
$ws.Range("B48").Value = 2024
$ws.Range("C48").Value = -0.4137309550271362
$ws.Range("D48").Value = 2025
$ws.Range("E48").Value = 0.6619818620803297

$ws.Range("A49").Value = 45618
$ws.Range("B49").Value = 2024
$ws.Range("C49").ClearContents()
$ws.Range("D49").Value = 2025
$ws.Range("E49").Value = 0.8860470190541037

$ws.Range("A50").Value = 45713
$ws.Range("B50").Value = 2025
$ws.Range("C50").Value = 0.3127680745294459
$ws.Range("D50").Value = 2026
$ws.Range("E50").Value = -0.5546801915590427

$ws.Range("A51").Value = 45800
$ws.Range("B51").Value = 2025
$ws.Range("C51").Value = 0.2094327661663842
$ws.Range("D51").Value = 2026
$ws.Range("E51").Value = -0.6242159253788016

$ws.Range("A52").Value = 45891
$ws.Range("B52").Value = 2025
$ws.Range("C52").Value = 0.2267356977060819
$ws.Range("D52").Value = 2026
$ws.Range("E52").Value = -0.4189713395563288
